$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1757.8541  # ALC!H17: 1757.0426 -> 1757.8541
$ws.Cells.Item(17, 10).Value = 1904.6578  # ALC!J17: 1907.5946 -> 1904.6578
$ws.Cells.Item(17, 12).Value = 5713.9734  # ALC!L17: 5722.783799999999 -> 5713.9734
$ws.Cells.Item(17, 14).Value = -6049.9734  # ALC!N17: -6058.783799999999 -> -6049.9734

$ws.Cells.Item(38, 8).Value = 103.55556  # ALC!H38: 760 -> 103.55556
$ws.Cells.Item(38, 9).Value = 103.55556  # ALC!I38: 116.375 -> 103.55556
$ws.Cells.Item(38, 10).Value = 0  # ALC!J38: 5909 -> 0
$ws.Cells.Item(38, 11).Value = 310.66668  # ALC!K38: 349.125 -> 310.66668
$ws.Cells.Item(38, 12).Value = 0  # ALC!L38: 17727 -> 0
$ws.Cells.Item(38, 13).Value = 61.33332000000001  # ALC!M38: 22.875 -> 61.33332000000001
$ws.Cells.Item(38, 14).ClearContents()  # ALC!N38 was -18471

$ws.Cells.Item(39, 8).Value = 541.75  # ALC!H39: 1123.1 -> 541.75
$ws.Cells.Item(39, 9).Value = 110.375  # ALC!I39: 810.2857 -> 110.375
$ws.Cells.Item(39, 10).Value = 1404.5  # ALC!J39: 1853 -> 1404.5
$ws.Cells.Item(39, 11).Value = 331.125  # ALC!K39: 2430.8571 -> 331.125
$ws.Cells.Item(39, 12).Value = 4213.5  # ALC!L39: 5559 -> 4213.5
$ws.Cells.Item(39, 13).Value = -35.125  # ALC!M39: -2134.8571 -> -35.125
$ws.Cells.Item(39, 14).Value = -4805.5  # ALC!N39: -6151 -> -4805.5

$ws.Cells.Item(74, 8).Value = 5951.0713  # ALC!H74: 6193.5386 -> 5951.0713
$ws.Cells.Item(74, 9).Value = 3119.8  # ALC!I74: 3200 -> 3119.8
$ws.Cells.Item(74, 11).Value = 3119.8  # ALC!K74: 3200 -> 3119.8
$ws.Cells.Item(74, 13).Value = -2183.8  # ALC!M74: -2264 -> -2183.8

$ws.Cells.Item(77, 8).Value = 5951.0713  # ALC!H77: 6193.5386 -> 5951.0713
$ws.Cells.Item(77, 9).Value = 3119.8  # ALC!I77: 3200 -> 3119.8
$ws.Cells.Item(77, 11).Value = 15599  # ALC!K77: 16000 -> 15599
$ws.Cells.Item(77, 13).Value = -10919  # ALC!M77: -11320 -> -10919

$ws.Cells.Item(86, 8).Value = 4301.067  # ALC!H86: 4565.4287 -> 4301.067
$ws.Cells.Item(86, 9).Value = 1775  # ALC!I86: 2166.6667 -> 1775
$ws.Cells.Item(86, 11).Value = 1775  # ALC!K86: 2166.6667 -> 1775
$ws.Cells.Item(86, 13).Value = -652  # ALC!M86: -1043.6667 -> -652

$ws.Cells.Item(89, 8).Value = 4301.067  # ALC!H89: 4565.4287 -> 4301.067
$ws.Cells.Item(89, 9).Value = 1775  # ALC!I89: 2166.6667 -> 1775
$ws.Cells.Item(89, 11).Value = 8875  # ALC!K89: 10833.3335 -> 8875
$ws.Cells.Item(89, 13).Value = -3259  # ALC!M89: -5217.333500000001 -> -3259

$ws.Cells.Item(106, 8).Value = 2321.818  # ALC!H106: 2245 -> 2321.818
$ws.Cells.Item(106, 9).Value = 1508  # ALC!I106: 1490 -> 1508
$ws.Cells.Item(106, 11).Value = 1508  # ALC!K106: 1490 -> 1508
$ws.Cells.Item(106, 13).Value = -877  # ALC!M106: -859 -> -877

$ws.Cells.Item(132, 8).Value = 28035.023  # ALC!H132: 29378.94 -> 28035.023
$ws.Cells.Item(132, 9).Value = 31336.377  # ALC!I132: 33045.29 -> 31336.377
$ws.Cells.Item(132, 11).Value = 94009.13099999999  # ALC!K132: 99135.87 -> 94009.13099999999
$ws.Cells.Item(132, 13).Value = -91479.13099999999  # ALC!M132: -96605.87 -> -91479.13099999999

$ws.Cells.Item(137, 8).Value = 1049601.2  # ALC!H137: 1117253.1 -> 1049601.2
$ws.Cells.Item(137, 9).Value = 1507  # ALC!I137: 1543.5714 -> 1507
$ws.Cells.Item(137, 10).Value = 1923013.2  # ALC!J137: 2036072.9 -> 1923013.2
$ws.Cells.Item(137, 11).Value = 4521  # ALC!K137: 4630.7142 -> 4521
$ws.Cells.Item(137, 12).Value = 5769039.6  # ALC!L137: 6108218.699999999 -> 5769039.6
$ws.Cells.Item(137, 13).Value = -1971  # ALC!M137: -2080.7142 -> -1971
$ws.Cells.Item(137, 14).Value = -5774139.6  # ALC!N137: -6113318.699999999 -> -5774139.6

$ws.Cells.Item(138, 8).Value = 2960.3635  # ALC!H138: 3077.2188 -> 2960.3635
$ws.Cells.Item(138, 9).Value = 2466.2666  # ALC!I138: 2555.2856 -> 2466.2666
$ws.Cells.Item(138, 10).Value = 3372.111  # ALC!J138: 3483.1667 -> 3372.111
$ws.Cells.Item(138, 11).Value = 7398.7998  # ALC!K138: 7665.8568 -> 7398.7998
$ws.Cells.Item(138, 12).Value = 10116.333  # ALC!L138: 10449.5001 -> 10116.333
$ws.Cells.Item(138, 13).Value = -2258.7998  # ALC!M138: -2525.8568 -> -2258.7998
$ws.Cells.Item(138, 14).Value = -20396.333  # ALC!N138: -20729.5001 -> -20396.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2598  # ARM!H2: 2075.3333 -> 2598
$ws.Cells.Item(2, 9).Value = 2664.3333  # ARM!I2: 2040.8462 -> 2664.3333
$ws.Cells.Item(2, 11).Value = 2664.3333  # ARM!K2: 2040.8462 -> 2664.3333
$ws.Cells.Item(2, 13).Value = -2551.3333  # ARM!M2: -1927.8462 -> -2551.3333

$ws.Cells.Item(32, 8).Value = 4247.7  # ARM!H32: 5666.03 -> 4247.7
$ws.Cells.Item(32, 9).Value = 1774  # ARM!I32: 2355.4268 -> 1774
$ws.Cells.Item(32, 10).Value = 12081.083  # ARM!J32: 15597.84 -> 12081.083
$ws.Cells.Item(32, 11).Value = 1774  # ARM!K32: 2355.4268 -> 1774
$ws.Cells.Item(32, 12).Value = 12081.083  # ARM!L32: 15597.84 -> 12081.083
$ws.Cells.Item(32, 13).Value = -1487  # ARM!M32: -2068.4268 -> -1487
$ws.Cells.Item(32, 14).Value = -12655.083  # ARM!N32: -16171.84 -> -12655.083

$ws.Cells.Item(74, 8).Value = 2634.5225  # ARM!H74: 2671.4243 -> 2634.5225
$ws.Cells.Item(74, 9).Value = 2488.65  # ARM!I74: 2547.359 -> 2488.65
$ws.Cells.Item(74, 11).Value = 2488.65  # ARM!K74: 2547.359 -> 2488.65
$ws.Cells.Item(74, 13).Value = -1614.65  # ARM!M74: -1673.359 -> -1614.65

$ws.Cells.Item(77, 8).Value = 2634.5225  # ARM!H77: 2671.4243 -> 2634.5225
$ws.Cells.Item(77, 9).Value = 2488.65  # ARM!I77: 2547.359 -> 2488.65
$ws.Cells.Item(77, 11).Value = 12443.25  # ARM!K77: 12736.795 -> 12443.25
$ws.Cells.Item(77, 13).Value = -8075.25  # ARM!M77: -8368.795 -> -8075.25

$ws.Cells.Item(116, 8).Value = 2598  # ARM!H116: 2075.3333 -> 2598
$ws.Cells.Item(116, 9).Value = 2664.3333  # ARM!I116: 2040.8462 -> 2664.3333
$ws.Cells.Item(116, 11).Value = 2664.3333  # ARM!K116: 2040.8462 -> 2664.3333
$ws.Cells.Item(116, 13).Value = -370.3332999999998  # ARM!M116: 253.1538 -> -370.3332999999998

$ws.Cells.Item(132, 8).Value = 2418.8772  # ARM!H132: 2411.4656 -> 2418.8772
$ws.Cells.Item(132, 9).Value = 1933.2766  # ARM!I132: 1934.4375 -> 1933.2766
$ws.Cells.Item(132, 11).Value = 5799.8298  # ARM!K132: 5803.3125 -> 5799.8298
$ws.Cells.Item(132, 13).Value = -3269.8298  # ARM!M132: -3273.3125 -> -3269.8298

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2598  # BSM!H3: 2075.3333 -> 2598
$ws.Cells.Item(3, 9).Value = 2664.3333  # BSM!I3: 2040.8462 -> 2664.3333
$ws.Cells.Item(3, 11).Value = 2664.3333  # BSM!K3: 2040.8462 -> 2664.3333
$ws.Cells.Item(3, 13).Value = -2550.3333  # BSM!M3: -1926.8462 -> -2550.3333

$ws.Cells.Item(20, 8).Value = 1856.8182  # BSM!H20: 1888.6666 -> 1856.8182
$ws.Cells.Item(20, 9).Value = 1824.2632  # BSM!I20: 1859.6111 -> 1824.2632
$ws.Cells.Item(20, 11).Value = 1824.2632  # BSM!K20: 1859.6111 -> 1824.2632
$ws.Cells.Item(20, 13).Value = -1577.2632  # BSM!M20: -1612.6111 -> -1577.2632

$ws.Cells.Item(105, 8).Value = 1853.9259  # BSM!H105: 1922.3077 -> 1853.9259
$ws.Cells.Item(105, 9).Value = 1789.2084  # BSM!I105: 1936.8572 -> 1789.2084
$ws.Cells.Item(105, 10).Value = 2371.6667  # BSM!J105: 1861.2 -> 2371.6667
$ws.Cells.Item(105, 11).Value = 1789.2084  # BSM!K105: 1936.8572 -> 1789.2084
$ws.Cells.Item(105, 12).Value = 2371.6667  # BSM!L105: 1861.2 -> 2371.6667
$ws.Cells.Item(105, 13).Value = -42.20839999999998  # BSM!M105: -189.8571999999999 -> -42.20839999999998
$ws.Cells.Item(105, 14).Value = -5865.6667  # BSM!N105: -5355.2 -> -5865.6667

$ws.Cells.Item(107, 8).Value = 2191.9285  # BSM!H107: 1248.3667 -> 2191.9285
$ws.Cells.Item(107, 9).Value = 2191.9285  # BSM!I107: 1281.069 -> 2191.9285
$ws.Cells.Item(107, 10).Value = 0  # BSM!J107: 300 -> 0
$ws.Cells.Item(107, 11).Value = 2191.9285  # BSM!K107: 1281.069 -> 2191.9285
$ws.Cells.Item(107, 12).Value = 0  # BSM!L107: 300 -> 0
$ws.Cells.Item(107, 13).Value = -271.9285  # BSM!M107: 638.931 -> -271.9285
$ws.Cells.Item(107, 14).ClearContents()  # BSM!N107 was -4140

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 41582.832  # CRP!H50: 38899.6 -> 41582.832
$ws.Cells.Item(50, 9).Value = 6749.5  # CRP!I50: 8749.5 -> 6749.5
$ws.Cells.Item(50, 10).Value = 58999.5  # CRP!J50: 58999.668 -> 58999.5
$ws.Cells.Item(50, 11).Value = 6749.5  # CRP!K50: 8749.5 -> 6749.5
$ws.Cells.Item(50, 12).Value = 58999.5  # CRP!L50: 58999.668 -> 58999.5
$ws.Cells.Item(50, 13).Value = -6124.5  # CRP!M50: -8124.5 -> -6124.5
$ws.Cells.Item(50, 14).Value = -60249.5  # CRP!N50: -60249.668 -> -60249.5

$ws.Cells.Item(132, 8).Value = 1433.7941  # CRP!H132: 1417.0857 -> 1433.7941
$ws.Cells.Item(132, 9).Value = 1517.9259  # CRP!I132: 1494.0358 -> 1517.9259
$ws.Cells.Item(132, 11).Value = 4553.7777  # CRP!K132: 4482.107400000001 -> 4553.7777
$ws.Cells.Item(132, 13).Value = -2023.7777  # CRP!M132: -1952.107400000001 -> -2023.7777

$ws.Cells.Item(140, 8).Value = 92361.25  # CRP!H140: 92307.09 -> 92361.25
$ws.Cells.Item(140, 10).Value = 92361.25  # CRP!J140: 92307.09 -> 92361.25
$ws.Cells.Item(140, 12).Value = 92361.25  # CRP!L140: 92307.09 -> 92361.25
$ws.Cells.Item(140, 14).Value = -102721.25  # CRP!N140: -102667.09 -> -102721.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(70, 8).Value = 3556.5557  # CUL!H70: 3001.75 -> 3556.5557

$ws.Cells.Item(73, 8).Value = 3556.5557  # CUL!H73: 3001.75 -> 3556.5557

$ws.Cells.Item(105, 8).Value = 13004.2  # CUL!H105: 11807.1 -> 13004.2
$ws.Cells.Item(105, 10).Value = 18340.834  # CUL!J105: 16345.667 -> 18340.834
$ws.Cells.Item(105, 12).Value = 55022.50199999999  # CUL!L105: 49037.001 -> 55022.50199999999
$ws.Cells.Item(105, 14).Value = -60264.50199999999  # CUL!N105: -54279.001 -> -60264.50199999999

$ws.Cells.Item(126, 8).Value = 5799.875  # CUL!H126: 7181.5 -> 5799.875
$ws.Cells.Item(126, 9).Value = 6223.8  # CUL!I126: 7181.5 -> 6223.8
$ws.Cells.Item(126, 10).Value = 5093.3335  # CUL!J126: 0 -> 5093.3335
$ws.Cells.Item(126, 11).Value = 18671.4  # CUL!K126: 21544.5 -> 18671.4
$ws.Cells.Item(126, 12).Value = 15280.0005  # CUL!L126: 0 -> 15280.0005
$ws.Cells.Item(126, 13).Value = -13731.4  # CUL!M126: -16604.5 -> -13731.4
$ws.Cells.Item(126, 14).Value = -25160.0005  # CUL!N126: None -> -25160.0005

$ws.Cells.Item(131, 8).Value = 11370.917  # CUL!H131: 11782.695 -> 11370.917
$ws.Cells.Item(131, 10).Value = 12361.454  # CUL!J131: 12859.619 -> 12361.454
$ws.Cells.Item(131, 12).Value = 37084.362  # CUL!L131: 38578.857 -> 37084.362
$ws.Cells.Item(131, 14).Value = -47164.362  # CUL!N131: -48658.857 -> -47164.362

$ws.Cells.Item(133, 8).Value = 3390.5557  # CUL!H133: 3327.9473 -> 3390.5557
$ws.Cells.Item(133, 9).Value = 3390.5557  # CUL!I133: 3327.9473 -> 3390.5557
$ws.Cells.Item(133, 11).Value = 10171.6671  # CUL!K133: 9983.841899999999 -> 10171.6671
$ws.Cells.Item(133, 13).Value = -5111.667099999999  # CUL!M133: -4923.841899999999 -> -5111.667099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5956.143  # GSM!H70: 5965.8335 -> 5956.143
$ws.Cells.Item(70, 9).Value = 5899  # GSM!I70: 5900 -> 5899
$ws.Cells.Item(70, 11).Value = 5899  # GSM!K70: 5900 -> 5899
$ws.Cells.Item(70, 13).Value = -5629  # GSM!M70: -5630 -> -5629

$ws.Cells.Item(73, 8).Value = 5956.143  # GSM!H73: 5965.8335 -> 5956.143
$ws.Cells.Item(73, 9).Value = 5899  # GSM!I73: 5900 -> 5899
$ws.Cells.Item(73, 11).Value = 5899  # GSM!K73: 5900 -> 5899
$ws.Cells.Item(73, 13).Value = -4963  # GSM!M73: -4964 -> -4963

$ws.Cells.Item(122, 8).Value = 4349.0645  # GSM!H122: 4434.467 -> 4349.0645
$ws.Cells.Item(122, 9).Value = 4504.5186  # GSM!I122: 4609.0386 -> 4504.5186
$ws.Cells.Item(122, 11).Value = 13513.5558  # GSM!K122: 13827.1158 -> 13513.5558
$ws.Cells.Item(122, 13).Value = -11063.5558  # GSM!M122: -11377.1158 -> -11063.5558

$ws.Cells.Item(132, 8).Value = 671563.75  # GSM!H132: 355813.1 -> 671563.75
$ws.Cells.Item(132, 9).Value = 928726.7  # GSM!I132: 389959.47 -> 928726.7
$ws.Cells.Item(132, 10).Value = 2940  # GSM!J132: 2967 -> 2940
$ws.Cells.Item(132, 11).Value = 2786180.1  # GSM!K132: 1169878.41 -> 2786180.1
$ws.Cells.Item(132, 12).Value = 8820  # GSM!L132: 8901 -> 8820
$ws.Cells.Item(132, 13).Value = -2783650.1  # GSM!M132: -1167348.41 -> -2783650.1
$ws.Cells.Item(132, 14).Value = -13880  # GSM!N132: -13961 -> -13880

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1557.0625  # LTW!H61: 1707.8572 -> 1557.0625
$ws.Cells.Item(61, 9).Value = 1539.5555  # LTW!I61: 1836.1428 -> 1539.5555
$ws.Cells.Item(61, 11).Value = 1539.5555  # LTW!K61: 1836.1428 -> 1539.5555
$ws.Cells.Item(61, 13).Value = -1337.5555  # LTW!M61: -1634.1428 -> -1337.5555

$ws.Cells.Item(113, 8).Value = 1557.0625  # LTW!H113: 1707.8572 -> 1557.0625
$ws.Cells.Item(113, 9).Value = 1539.5555  # LTW!I113: 1836.1428 -> 1539.5555
$ws.Cells.Item(113, 11).Value = 1539.5555  # LTW!K113: 1836.1428 -> 1539.5555
$ws.Cells.Item(113, 13).Value = 630.4445000000001  # LTW!M113: 333.8571999999999 -> 630.4445000000001

$ws.Cells.Item(122, 8).Value = 39825.965  # LTW!H122: 41286.332 -> 39825.965
$ws.Cells.Item(122, 9).Value = 3394  # LTW!I122: 3704 -> 3394
$ws.Cells.Item(122, 10).Value = 54398.75  # LTW!J122: 57110.473 -> 54398.75
$ws.Cells.Item(122, 11).Value = 10182  # LTW!K122: 11112 -> 10182
$ws.Cells.Item(122, 12).Value = 163196.25  # LTW!L122: 171331.419 -> 163196.25
$ws.Cells.Item(122, 13).Value = -7732  # LTW!M122: -8662 -> -7732
$ws.Cells.Item(122, 14).Value = -168096.25  # LTW!N122: -176231.419 -> -168096.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 1500  # WVR!H18: 0 -> 1500
$ws.Cells.Item(18, 9).Value = 1500  # WVR!I18: 0 -> 1500
$ws.Cells.Item(18, 11).Value = 1500  # WVR!K18: 0 -> 1500
$ws.Cells.Item(18, 13).Value = -1327  # WVR!M18: None -> -1327

$ws.Cells.Item(81, 8).Value = 933  # WVR!H81: 999.5 -> 933
$ws.Cells.Item(81, 9).Value = 933  # WVR!I81: 999.5 -> 933
$ws.Cells.Item(81, 11).Value = 1866  # WVR!K81: 1999 -> 1866
$ws.Cells.Item(81, 13).Value = -805  # WVR!M81: -938 -> -805

$ws.Cells.Item(84, 8).Value = 933  # WVR!H84: 999.5 -> 933
$ws.Cells.Item(84, 9).Value = 933  # WVR!I84: 999.5 -> 933
$ws.Cells.Item(84, 11).Value = 9330  # WVR!K84: 9995 -> 9330
$ws.Cells.Item(84, 13).Value = -4026  # WVR!M84: -4691 -> -4026

$ws.Cells.Item(126, 8).Value = 4799.2  # WVR!H126: 3445.2727 -> 4799.2
$ws.Cells.Item(126, 9).Value = 6000  # WVR!I126: 2843.1428 -> 6000
$ws.Cells.Item(126, 11).Value = 18000  # WVR!K126: 8529.428400000001 -> 18000
$ws.Cells.Item(126, 13).Value = -15530  # WVR!M126: -6059.428400000001 -> -15530

$ws.Cells.Item(132, 8).Value = 1649797.2  # WVR!H132: 1732282.2 -> 1649797.2
$ws.Cells.Item(132, 9).Value = 3144049.5  # WVR!I132: 3842494.2 -> 3144049.5
$ws.Cells.Item(132, 10).Value = 6119.7  # WVR!J132: 5745.1816 -> 6119.7
$ws.Cells.Item(132, 11).Value = 9432148.5  # WVR!K132: 11527482.6 -> 9432148.5
$ws.Cells.Item(132, 12).Value = 18359.1  # WVR!L132: 17235.5448 -> 18359.1
$ws.Cells.Item(132, 13).Value = -9429618.5  # WVR!M132: -11524952.6 -> -9429618.5
$ws.Cells.Item(132, 14).Value = -23419.1  # WVR!N132: -22295.5448 -> -23419.1
